# Regenerate merged AHB files
# - rename the _old/_new header-string families to _FV2310/_FV2404
# - wrap the data range in an Excel Table (with AutoFilter)
# - freeze the header row

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename header row (A1:J1 = _old -> _FV2310, L1:U1 = _new -> _FV2404; K1 "diff" stays) ---
$ws.Range("A1").Value = "Segmentname_FV2310"
$ws.Range("B1").Value = "Segmentgruppe_FV2310"
$ws.Range("C1").Value = "Segment_FV2310"
$ws.Range("D1").Value = "Datenelement_FV2310"
$ws.Range("E1").Value = "Segment ID_FV2310"
$ws.Range("F1").Value = "Code_FV2310"
$ws.Range("G1").Value = "Qualifier_FV2310"
$ws.Range("H1").Value = "Beschreibung_FV2310"
$ws.Range("I1").Value = "Bedingungsausdruck_FV2310"
$ws.Range("J1").Value = "Bedingung_FV2310"

$ws.Range("L1").Value = "Segmentname_FV2404"
$ws.Range("M1").Value = "Segmentgruppe_FV2404"
$ws.Range("N1").Value = "Segment_FV2404"
$ws.Range("O1").Value = "Datenelement_FV2404"
$ws.Range("P1").Value = "Segment ID_FV2404"
$ws.Range("Q1").Value = "Code_FV2404"
$ws.Range("R1").Value = "Qualifier_FV2404"
$ws.Range("S1").Value = "Beschreibung_FV2404"
$ws.Range("T1").Value = "Bedingungsausdruck_FV2404"
$ws.Range("U1").Value = "Bedingung_FV2404"

# --- 2. Turn the used range into a native Excel Table with an AutoFilter ---
$tbl = $ws.ListObjects.Add(1, $ws.Range("A1:U84"), $null, 1)
$tbl.Name = "Table1"
$tbl.TableStyle = "TableStyleMedium9"

# --- 3. Freeze the header row (ySplit=1, top-left cell of the scrolling pane = A2) ---
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

Write-Output "edit complete"
